$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold text (not numeric) values, e.g. "567.4" rather than 567.4,
# so format them as Text first to force the COM layer to keep the literal string
# instead of auto-converting the numeric-looking text into a number.
$rng = $ws.Range("C2:G4")
$rng.NumberFormat = "@"

$ws.Range("C2").Value = "567.4"
$ws.Range("D2").Value = "562.3"
$ws.Range("E2").Value = "562.8"
$ws.Range("F2").Value = "564.1"
$ws.Range("G2").Value = "563.9"

$ws.Range("C3").Value = "357.6"
$ws.Range("D3").Value = "355.1"
$ws.Range("E3").Value = "359.6"
$ws.Range("F3").Value = "358.0"
$ws.Range("G3").Value = "361.0"

$ws.Range("C4").Value = "196.2"
$ws.Range("D4").Value = "195.9"
$ws.Range("E4").Value = "199.0"
$ws.Range("F4").Value = "198.8"
$ws.Range("G4").Value = "199.9"
